$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.736.97'
$ws.Range('E2').Value = '  -2.90%  '

$ws.Range('D3').Value = '2.613.18'
$ws.Range('E3').Value = '  -1.65%  '

$ws.Range('E4').Value = '  +0.03%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '575.98'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -3.51%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '156.46'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.04%  '

$ws.Range('E7').Value = '  +0.05%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.622'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -3.80%  '

$ws.Range('E9').Value = '  -4.98%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.82'
$ws.Range('D10').Style = 'Normal'

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.382'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -4.27%  '

$ws.Range('E12').Value = '  -0.17%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '28.22'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.95%  '

$ws.Range('D14').Value = '3.088.93'
$ws.Range('E14').Value = '  -1.26%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000182'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -6.62%  '

$ws.Range('D16').Value = '63.556.69'
$ws.Range('E16').Value = '  -2.96%  '

$ws.Range('D17').Value = '2.612.24'
$ws.Range('E17').Value = '  -2.73%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '12.09'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -3.44%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.69'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.74%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.56'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -4.28%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '342.63'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.92%  '

$ws.Range('E22').Value = '  -0.22%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '67.44'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.94%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.76'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -1.10%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.0000109'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.52%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '592.54'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +5.65%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.22'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.04%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.58'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.59%  '

$ws.Range('E29').Value = '  -0.12%  '

$ws.Range('B30').Value = 'Aptos'
$ws.Range('C30').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.93'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.46%  '

$ws.Range('B31').Value = 'Kaspa'
$ws.Range('C31').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.160'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.19%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.06'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -2.68%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.76'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.80%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '6.61'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.32%  '

$ws.Range('E35').Value = '  -1.59%  '

$ws.Range('E36').Value = '  -2.82%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '19.75'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -3.16%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.999'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -0.03%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '154.52'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.45%  '

$ws.Range('E40').Value = '  -3.01%  '

$ws.Range('E41').Value = '  -0.01%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '41.53'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.66%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.43'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +7.55%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '155.89'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -2.76%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '3.91'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -3.68%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '23.36'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +3.60%  '

$ws.Range('E47').Value = '  -1.95%  '

$ws.Range('E48').Value = '  +0.05%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.628'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -1.54%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0247'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.73%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '18.94'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.40%  '
